# Fix spelling errors in the presentation:
#  - Slide 21: "Test code coverage and other matrices" -> split off "other "
#    into its own run and fix "matrices" -> "metrics"
#  - Slide 24: "Question?" -> "Questions?"

$p = $ppt.ActivePresentation

# --- Slide 21: "... other matrices" -> "... other metrics" (split runs) ---
$s21 = $p.Slides.Item(21)
$shp21 = $s21.Shapes.Item(2)
$tr21 = $shp21.TextFrame.TextRange

# Full text (paragraphs joined by CR) is:
#   "MVC Framework itself\rDRY\rOCP\rTest code coverage and other matrices"
# The run "Test code coverage and other " starts at character 30 (1-based)
# and is 29 characters long; "other " occupies the last 6 of those (chars
# 53-58). The following run "matrices" starts at character 59.

# Split "other " out of the "Test code coverage and other " run so it
# becomes its own run.
$otherRun = $tr21.Characters(53, 6)
$otherRun.Text = "other "

# Fix the spelling of "matrices" -> "metrics" in the last run.
$metricsRun = $tr21.Characters(59, 8)
$metricsRun.Text = "metrics"

# --- Slide 24: "Question?" -> "Questions?" ---
$s24 = $p.Slides.Item(24)
$shp24 = $s24.Shapes.Item(1)
$shp24.TextFrame.TextRange.Text = "Questions?"
